$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "59.311.03"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.65%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.634.18"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.76%  "

$ws.Range("E4").Value = "  +0.09%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "515.83"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "146.19"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.98%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "

$ws.Range("E8").Value = "  +0.10%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.656.70"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.50%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.46"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.03%  "

$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("E13").Value = "  -1.41%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "3.101.54"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.01%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "59.275.45"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.73%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.21"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.74%  "

$ws.Range("E17").Value = "  -0.19%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "2.645.95"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.37%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.60"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.77%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "345.06"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.15%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "10.47"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "6.15"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.64%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +1.05%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "60.79"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.95%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.423"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.73%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.795.48"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +2.51%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.20%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.160"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.84%  "

$ws.Range("E29").Value = "  +1.06%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.23"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.13%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +7.89%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "18.97"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.02%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.58"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.06%  "

$ws.Range("E35").Value = "  +15.71%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "149.10"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.52%  "

$ws.Range("E37").Value = "  +1.22%  "

$ws.Range("E38").Value = "  +1.38%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.867"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "36.56"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.12%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.72"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.46%  "

$ws.Range("E42").Value = "  -1.36%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "281.64"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.617"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.40%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.0993"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -1.39%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.994"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -0.45%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "19.66"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("E48").Value = "  -0.80%  "

$ws.Range("E49").Value = "  -0.52%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.74"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.83%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "10.25"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.39%  "
